# Refresh the cryptos price list (Coin/Price/Volume(1h)) with the latest
# snapshot. Column D ("Price") and E ("Volume(1h)") cells are stored as
# plain text in the sheet (e.g. "71.026.20", "  +2.72%  "), so decimal-
# looking values are written with a leading apostrophe to force Excel to
# keep them as literal text instead of silently parsing/rounding them as
# numbers (which would drop significant trailing/leading zeros, e.g.
# "11.20" -> 11.2 or "0.300" -> 0.3).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.026.20"
$ws.Range("E2").Value = "  +2.72%  "

$ws.Range("D3").Value = "3.790.26"
$ws.Range("E3").Value = "  +0.48%  "

$ws.Range("E4").Value = "  +0.21%  "

$ws.Range("D5").Value = "'702.35"
$ws.Range("E5").Value = "  +11.34%  "

$ws.Range("D6").Value = "'172.94"

$ws.Range("D7").Value = "3.788.07"
$ws.Range("E7").Value = "  +0.45%  "

$ws.Range("E8").Value = "  -0.08%  "

$ws.Range("E9").Value = "  +0.80%  "

$ws.Range("E10").Value = "  +2.54%  "

$ws.Range("D11").Value = "'7.47"
$ws.Range("E11").Value = "  +10.28%  "

$ws.Range("E12").Value = "  +0.33%  "

$ws.Range("D13").Value = "'0.0000256"
$ws.Range("E13").Value = "  +6.72%  "

$ws.Range("D14").Value = "'36.12"
$ws.Range("E14").Value = "  +3.31%  "

$ws.Range("D15").Value = "4.428.28"
$ws.Range("E15").Value = "  +0.56%  "

$ws.Range("D16").Value = "3.793.58"
$ws.Range("E16").Value = "  -0.20%  "

$ws.Range("D17").Value = "70.990.98"
$ws.Range("E17").Value = "  +2.72%  "

$ws.Range("D18").Value = "'17.88"
$ws.Range("E18").Value = "  +1.44%  "

$ws.Range("D19").Value = "'7.19"
$ws.Range("E19").Value = "  +2.48%  "

$ws.Range("E20").Value = "  +0.80%  "

$ws.Range("D21").Value = "'11.20"
$ws.Range("E21").Value = "  +17.61%  "

$ws.Range("D22").Value = "'482.63"
$ws.Range("E22").Value = "  +4.42%  "

$ws.Range("E23").Value = "  +1.18%  "

$ws.Range("D24").Value = "'83.85"
$ws.Range("E24").Value = "  +2.21%  "

$ws.Range("E25").Value = "  +0.21%  "

$ws.Range("E26").Value = "  +2.03%  "

$ws.Range("B27").Value = "Fetch.AI"
$ws.Range("C27").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D27").Value = "'2.19"
$ws.Range("E27").Value = "  +2.18%  "

$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").Value = "'10.53"
$ws.Range("E28").Value = "  +4.44%  "

$ws.Range("D29").Value = "3.939.87"
$ws.Range("E29").Value = "  +0.53%  "

$ws.Range("D30").Value = "'0.999"
$ws.Range("E30").Value = "  -0.11%  "

$ws.Range("D31").Value = "'3.11"
$ws.Range("E31").Value = "  +15.74%  "

$ws.Range("E32").Value = "  +1.09%  "

$ws.Range("D33").Value = "'7.55"
$ws.Range("E33").Value = "  +6.96%  "

$ws.Range("D34").Value = "'29.52"
$ws.Range("E34").Value = "  +3.67%  "

$ws.Range("D35").Value = "'0.178"
$ws.Range("E35").Value = "  -1.12%  "

$ws.Range("D36").Value = "'9.19"
$ws.Range("E36").Value = "  +2.67%  "

$ws.Range("E37").Value = "  +0.07%  "

$ws.Range("D38").Value = "3.739.25"
$ws.Range("E38").Value = "  +0.45%  "

$ws.Range("E39").Value = "  +2.06%  "

$ws.Range("D40").Value = "'3.45"
$ws.Range("E40").Value = "  +5.00%  "

$ws.Range("E41").Value = "  +3.06%  "

$ws.Range("D42").Value = "'2.23"
$ws.Range("E42").Value = "  +12.80%  "

$ws.Range("D43").Value = "'0.000325"
$ws.Range("E43").Value = "  +22.38%  "

$ws.Range("D44").Value = "'0.965"
$ws.Range("E44").Value = "  +0.25%  "

$ws.Range("E45").Value = "  +0.02%  "

$ws.Range("D47").Value = "'46.13"
$ws.Range("E47").Value = "  +6.88%  "

$ws.Range("D48").Value = "'161.29"
$ws.Range("E48").Value = "  +2.20%  "

$ws.Range("D49").Value = "'49.06"
$ws.Range("E49").Value = "  +4.34%  "

$ws.Range("E50").Value = "  -1.23%  "

$ws.Range("D51").Value = "'0.300"
$ws.Range("E51").Value = "  +1.69%  "
